$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,E,F,I,J,K,L,M,N (G stays 1, H is empty)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @{
    2 = @(1.02, 0.9999970461623656, 1.026513451475777, 1.003512603051345, 0.9980539377174651, 1.028140619206577, 1.005312350954267, 1.029335878222202, 1.006403886417997, 1.000962071010937, 1.005586813710726)
    3 = @(1.02, 1.002190696077171, 1.026832767717716, 1.005418250165748, 1.000925108801209, 1.028119731977348, 1.007127342732227, 1.029464047931005, 1.008108944823934, 1.003628619941815, 1.006228286426888)
    4 = @(1.02, 1.003602290855779, 1.027038990663314, 1.006644382404554, 1.002773591153557, 1.02810454699596, 1.008294031562883, 1.029545830079758, 1.009204923891012, 1.005344591826161, 1.006639578476929)
    5 = @(1.02, 1.004193894425286, 1.027125589626534, 1.007158223888107, 1.003548520839981, 1.028097762610924, 1.008782696506916, 1.029579933844558, 1.009663960881819, 1.006063784403245, 1.006811593515251)
    6 = @(1.02, 1.004293121253613, 1.027140124207022, 1.007244405997112, 1.003678509513614, 1.028096599963332, 1.008864640454659, 1.029585643700014, 1.009740935846062, 1.006184412705869, 1.006840423675391)
    7 = @(1.02, 1.003610203023348, 1.027040148187642, 1.006651254710766, 1.002783954247485, 1.028104457917983, 1.008300568190578, 1.029546286868134, 1.009211064246002, 1.005354210280489, 1.006641880437511)
    8 = @(1.019999999999999, 1.000740052902658, 1.026621446480098, 1.004158087951371, 0.9990262421512628, 1.02813390540539, 1.005927361266509, 1.029379431388346, 1.006981654404119, 1.001865241835927, 1.005804394172819)
    9 = @(1.02, 0.9956202764081897, 1.025880711243438, 0.9997098096606482, 0.9923299142558761, 1.028173049635261, 1.001684422044947, 1.029076656696341, 1.002995500289088, 0.9956418554368396, 1.004299033365596)
    10 = @(1.02, 0.9921622425845256, 1.025385058414673, 0.9967047873794241, 0.9878110444644356, 1.028190631074268, 0.9988121622532348, 1.02886901662911, 1.000296919877218, 0.9914381395819487, 1.003274668619613)
    11 = @(1.02, 0.9906535521406715, 1.025170031100778, 0.9953936374089645, 0.9858403664321082, 1.028196236062541, 0.9975575053047278, 1.028777753015332, 0.9991181050659224, 0.9896039559150833, 1.002825974291096)
    12 = @(1.02, 0.9900913939142384, 1.025090101753847, 0.9949050722402973, 0.9851061820500319, 1.028198017265125, 0.9970897735031714, 1.028743651841868, 0.9986786442991347, 0.9889204816658863, 1.002658519692851)
    13 = @(1.02, 0.9902120595381855, 1.025107249484422, 0.9950099418895325, 0.9852637674690576, 1.028197648785465, 0.9971901811635615, 1.028750975765279, 0.9987729831400927, 0.9890671890455449, 1.00269447529006)
    14 = @(1.02, 0.9906071202911284, 1.025163425298077, 0.9953532842849065, 0.9857797236426981, 1.028196389426981, 0.9975188773544146, 1.028774938309027, 0.9990818119466078, 0.9895475046040829, 1.002812148652134)
    15 = @(1.02, 0.9908502947745117, 1.025198029374264, 0.9955646225320298, 0.9860973288618624, 1.028195573671256, 0.9977211715879997, 1.028789675730796, 0.9992718785883367, 0.9898431519607224, 1.002884545935586)
    16 = @(1.02, 0.9922621254230287, 1.025399320714895, 0.9967915900936591, 0.9879415300188076, 1.028190216864643, 0.9988951946573905, 1.028875045105079, 1.00037493267123, 0.9915595675619321, 1.003304337241511)
    17 = @(1.02, 0.9931446518951077, 1.025525478374841, 0.9975585341355973, 0.989094545834366, 1.028186319765488, 0.9996286609749276, 1.02892823368536, 1.001064055618971, 0.992632438877268, 1.003566273402792)
    18 = @(1.02, 0.993658324721283, 1.025599024633332, 0.9980049224938602, 0.9897657377869169, 1.02818385285959, 1.000055426526088, 1.02895912692504, 1.001465017693937, 0.9932568864796915, 1.003718561842482)
    19 = @(1.02, 0.9938332907999199, 1.025624095254975, 0.9981569685386646, 0.9899943717553462, 1.028182978809425, 1.000200765487866, 1.028969638492465, 1.001601568606567, 0.9934695821076397, 1.003770404919007)
    20 = @(1.02, 0.9930500783226616, 1.025511946906639, 0.9974763476563392, 0.9889709776493107, 1.028186757925705, 0.9995500762588492, 1.028922540567341, 1.00099022216119, 0.9925174691326852, 1.003538221405817)
    21 = @(1.02, 0.9904908339126361, 1.025146884520153, 0.9952522215828185, 0.9856278484838727, 1.028196768571165, 0.9974221317641361, 1.028767887494854, 0.9989909140153368, 0.9894061243963045, 1.002777518719308)
    22 = @(1.02, 0.9888715021031147, 1.024917017289605, 0.9938448567504581, 0.9835131937094171, 1.028201323127522, 0.9960743677744542, 1.028669484195791, 0.9977246088080334, 0.9874372626241351, 1.002294657611214)
    23 = @(1.02, 0.9897309308844535, 1.025038905493158, 0.9945917945827535, 0.9846354452697469, 1.028199073251875, 0.9967897932096135, 1.028721759700692, 0.998396795041804, 0.9884822189254497, 1.002551071293827)
    24 = @(1.02, 0.9930928154140612, 1.025518061316701, 0.9975134871437007, 0.9890268169374212, 1.028186560538906, 0.9995855885604165, 1.028925113446604, 1.001023587384582, 0.9925694231479877, 1.003550898422906)
    25 = @(1.02, 0.9969515404821367, 1.026072541693807, 1.000866570495859, 0.9940703804778075, 1.028164434695646, 1.002788820867278, 1.029155958733625, 1.004033090616584, 0.9972601003021522, 1.004691803431707)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $vals[$i]
    }
}

Write-Output "Updated vm_pu values for case with 380 kV"
